# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet, and
# refreshes the Priority + Latest Handoff Datetime columns for the zh-cn
# (and Priority for de-de) rows that were just re-handed-off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the four
# rows that were re-generated.
$overview.Range("G4:G7").Value = "2016-08-22 11:55:03"

# zh-cn sheet: Priority moved from "low" to "ht", and the handoff xliff was
# regenerated with a new timestamp.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-22 11:54:55"

# de-de sheet: Priority moved from "low" to "ht" as well. Its "Latest
# Handoff Datetime" column happens to share the same underlying value as
# the Overview's "Latest HO Xliff Generate Date", so it picks up the same
# refreshed timestamp.
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-22 11:55:03"
